# DC COVID-19 Data workbook update — add the 4/25/2020 (serial 43946) data
# column to every sheet (one day after the previous last column), plus a
# stray single-cell row 99 addition on "Overal Stats" (AZ99 = 9).
#
# Layout of each sheet (before -> after):
#   Overal Stats (sheet 1)       : last data col AY -> AZ   (date header row 1)
#   Total Cases by Ward (sheet 2): last data col Z  -> AA   (date header row 2)
#   Total Cases by Race (sheet 3): last data col U  -> V    (date header row 2)
#   Lives Lost by Race (sheet 4) : last data col U  -> V    (date header row 1)
#   Lives Lost by Ward (sheet 5) : last data col G  -> H    (date header row 2)

$wb = $excel.ActiveWorkbook

function Set-NewColumnData {
    param(
        $ws,
        [string]$NewCol,
        [string]$PrevCol,
        [int]$HeaderRow,
        [double]$HeaderValue,
        [array]$Rows   # array of @(row, value) pairs for the data rows
    )

    # Header (date) cell: same date-number-format style as the previous
    # header cell in that row.
    $headerNew = $ws.Range($NewCol + $HeaderRow)
    $headerNew.Value = $HeaderValue
    $headerNew.NumberFormat = $ws.Range($PrevCol + $HeaderRow).NumberFormat

    foreach ($pair in $Rows) {
        $r = $pair[0]
        $v = $pair[1]
        $ws.Range($NewCol + $r).Value = $v
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "Overal Stats" -> new column AZ (date 43946), plus row 99
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overal Stats")

$ws1Rows = @(
    @(3,18068), @(4,3841), @(5,178), @(6,657),
    @(8,113), @(9,440), @(10,201), @(11,239),
    @(15,87), @(16,36), @(17,51), @(18,172), @(19,208), @(20,659),
    @(23,90), @(24,50), @(25,40), @(26,80), @(27,130), @(28,713),
    @(31,37), @(32,29), @(33,0), @(34,32), @(35,61), @(36,119), @(37,1),
    @(40,121), @(41,54), @(42,68), @(43,792), @(44,846), @(45,244), @(46,1),
    @(49,22), @(50,17), @(51,4), @(52,37), @(53,55), @(54,84), @(55,1),
    @(57,9), @(58,8), @(59,1), @(60,40), @(61,48), @(62,0), @(63,1),
    @(66,152), @(67,248), @(68,244), @(69,9),
    @(71,69), @(72,21), @(73,89), @(74,93),
    @(76,46), @(77,29), @(78,34), @(79,2), @(80,8)
)

Set-NewColumnData -ws $ws1 -NewCol "AZ" -PrevCol "AY" -HeaderRow 1 -HeaderValue 43946 -Rows $ws1Rows

# Lone new row far below the table
$ws1.Range("AZ99").Value = 9

# Extend the bestFit column-width formatting band (cols 40-51) to include
# the new column 52 (AZ), matching the previous columns' width.
$ws1.Range("AZ1").ColumnWidth = $ws1.Range("AY1").ColumnWidth

# ---------------------------------------------------------------------
# Sheet 2: "Total Cases by Ward" -> new column AA (date 43946)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Total Cases by Ward")

$ws2Rows = @(
    @(3,444), @(4,261), @(5,268), @(6,663), @(7,552),
    @(8,419), @(9,570), @(10,530), @(11,134)
)

Set-NewColumnData -ws $ws2 -NewCol "AA" -PrevCol "Z" -HeaderRow 2 -HeaderValue 43946 -Rows $ws2Rows

# ---------------------------------------------------------------------
# Sheet 3: "Total Cases by Race" -> new column V (date 43946)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Total Cases by Race")

$ws3Rows = @(
    @(4,3841), @(5,527), @(6,665), @(7,1891), @(8,63),
    @(9,11), @(10,8), @(11,652), @(12,24),
    @(14,799), @(15,708), @(16,2320), @(17,14)
)

Set-NewColumnData -ws $ws3 -NewCol "V" -PrevCol "U" -HeaderRow 2 -HeaderValue 43946 -Rows $ws3Rows

# ---------------------------------------------------------------------
# Sheet 4: "Lives Lost by Race" -> new column V (date 43946)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Lives Lost by Race")

$ws4Rows = @(
    @(3,178), @(4,4), @(5,142), @(6,11), @(7,20), @(8,1)
)

Set-NewColumnData -ws $ws4 -NewCol "V" -PrevCol "U" -HeaderRow 1 -HeaderValue 43946 -Rows $ws4Rows

# ---------------------------------------------------------------------
# Sheet 5: "Lives Lost by Ward" -> new column H (date 43946)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Lives Lost by Ward")

$ws5Rows = @(
    @(3,178), @(4,11), @(5,6), @(6,12), @(7,22),
    @(8,25), @(9,27), @(10,25), @(11,41), @(12,9), @(13,0)
)

Set-NewColumnData -ws $ws5 -NewCol "H" -PrevCol "G" -HeaderRow 2 -HeaderValue 43946 -Rows $ws5Rows

# ---------------------------------------------------------------------
# View-state touch-ups (selection changes recorded in the saved session)
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("R29").Select()

$ws3.Activate()
$ws3.Range("V2").Select()

$ws4.Activate()
$ws4.Range("V7").Select()

$ws5.Activate()
$ws5.Range("H14").Select()

# Leave the workbook on sheet 1 with the new last column selected, matching
# the original file's "tabSelected" + active-cell state on this sheet.
$ws1.Activate()
$ws1.Range("AY80").Select()
